$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '26.268.16'
$ws.Range("E2").Value = '  -3.40%  '
$ws.Range("D3").Value = '1.791.18'
$ws.Range("E3").Value = '  -3.11%  '
$ws.Range("E4").Value = '  +0.70%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '1.007'
$ws.Range("E5").Value = '  +0.61%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '306.97'
$ws.Range("E6").Value = '  -1.95%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4548'
$ws.Range("E7").Value = '  -1.87%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3609'
$ws.Range("E8").Value = '  -2.37%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.07061'
$ws.Range("E9").Value = '  -2.92%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.8698'
$ws.Range("E10").Value = '  -1.88%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.07764'
$ws.Range("E11").Value = '  -0.84%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '19.30'
$ws.Range("E12").Value = '  -3.14%  '
$ws.Range("D13").Value = '1.785.18'
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '5.261'
$ws.Range("E14").Value = '  -2.38%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '6.310'
$ws.Range("E15").Value = '  -3.05%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '84.68'
$ws.Range("E16").Value = '  -7.44%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '1.009'
$ws.Range("E17").Value = '  +0.69%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.000008484'
$ws.Range("E18").Value = '  -4.06%  '
$ws.Range("E19").Value = '  +0.66%  '
$ws.Range("D20").Value = '26.331.24'
$ws.Range("E20").Value = '  -3.27%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '14.09'
$ws.Range("E21").Value = '  -3.77%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '4.968'
$ws.Range("E22").Value = '  -1.77%  '
$ws.Range("D23").Value = '2.034.70'
$ws.Range("E23").Value = '  -4.06%  '
$ws.Range("E24").Value = '  -0.56%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '1.975'
$ws.Range("E25").Value = '  -3.02%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '152.23'
$ws.Range("E26").Value = '  +0.58%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '17.77'
$ws.Range("E27").Value = '  -3.40%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '2.027'
$ws.Range("E28").Value = '  +0.26%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '112.16'
$ws.Range("E29").Value = '  -2.99%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '4.822'
$ws.Range("E30").Value = '  -3.69%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.08642'
$ws.Range("E31").Value = '  -2.18%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '3.025'
$ws.Range("E32").Value = '  -3.38%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '4.432'
$ws.Range("E33").Value = '  -1.74%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.7124'
$ws.Range("E34").Value = '  -8.79%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.098'
$ws.Range("E35").Value = '  -4.26%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '2.627'
$ws.Range("E36").Value = '  -2.29%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '1.007'
$ws.Range("E37").Value = '  +0.68%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '1.075'
$ws.Range("E38").Value = '  -2.45%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.01934'
$ws.Range("E39").Value = '  -0.42%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.05078'
$ws.Range("E40").Value = '  -2.58%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '2.862'
$ws.Range("E41").Value = '  -3.07%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '6.875'
$ws.Range("E42").Value = '  -2.23%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.4909'
$ws.Range("E43").Value = '  -2.60%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.1518'
$ws.Range("E44").Value = '  -5.85%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '7.952'
$ws.Range("E45").Value = '  -6.58%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '1.008'
$ws.Range("E46").Value = '  +0.69%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.4559'
$ws.Range("E47").Value = '  -4.15%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '9.826'
$ws.Range("E48").Value = '  -5.48%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '99.78'
$ws.Range("E49").Value = '  -3.11%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '1.577'
$ws.Range("E50").Value = '  -3.44%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.05943'
$ws.Range("E51").Value = '  -4.09%  '
